$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HCL")

# Row 7
$ws.Range("F7").Value = 1432.15
$ws.Range("G7").Value = 1447.05
$ws.Range("H7").Value = 1420.1
$ws.Range("I7").Value = 1422.05
$ws.Range("J7").Value = 1439.95

# Row 9
$ws.Range("G9").Value = 1436
$ws.Range("H9").Value = 1421.15
$ws.Range("I9").Value = 1431.8

# Row 10
$ws.Range("G10").Value = 1446.85
$ws.Range("H10").Value = 1428.6
$ws.Range("I10").Value = 1439.8

# Row 11
$ws.Range("G11").Value = 1447.05
$ws.Range("H11").Value = 1437.5
$ws.Range("I11").Value = 1441.5

# Row 12
$ws.Range("G12").Value = 1445.9
$ws.Range("H12").Value = 1436.4
$ws.Range("I12").Value = 1442.05

# Row 13
$ws.Range("G13").Value = 1442.95
$ws.Range("H13").Value = 1437.05
$ws.Range("I13").Value = 1439.5

# Row 14
$ws.Range("G14").Value = 1441.85
$ws.Range("H14").Value = 1432.95
$ws.Range("I14").Value = 1435.4

# Row 15
$ws.Range("G15").Value = 1438
$ws.Range("H15").Value = 1435.05
$ws.Range("I15").Value = 1436.45

# Row 16
$ws.Range("G16").Value = 1440.05
$ws.Range("H16").Value = 1433.6
$ws.Range("I16").Value = 1434.95

# Row 17
$ws.Range("G17").Value = 1436.4
$ws.Range("H17").Value = 1428.55
$ws.Range("I17").Value = 1430.65

# Row 18
$ws.Range("G18").Value = 1434.1
$ws.Range("H18").Value = 1425.25
$ws.Range("I18").Value = 1427.3

# Row 19
$ws.Range("G19").Value = 1428.7
$ws.Range("H19").Value = 1423.6
$ws.Range("I19").Value = 1426

# Row 20
$ws.Range("G20").Value = 1427.95
$ws.Range("H20").Value = 1420.45
$ws.Range("I20").Value = 1420.5

# Row 21
$ws.Range("G21").Value = 1424.5
$ws.Range("H21").Value = 1420.1
$ws.Range("I21").Value = 1424.25
